$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.916.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +1.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3902"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07937"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.916.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.916"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.059"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.021"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06754"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001038"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.926.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.459"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.124.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.056"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.438"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.210"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.095"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5878"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1877"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.381"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.906"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06883"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  -1.34%  "
